# Updates to "Greece Super League 1" sheet: re-ordering of specific match
# rows (the underlying scraped data for these fixtures was resequenced).
# For each affected row, every column except A (the fixed positional
# index) is replaced with the values that belong to a different row.
# Column A keeps its original value for its row; the data columns B..AD
# move between rows as described below:
#   112 <- 113 (and 113 <- 112)
#   124 <- 125 (and 125 <- 124)
#   143 <- 144 (and 144 <- 143)
#   175 <- 180, 179 <- 175, 180 <- 179   (3-way rotation)
#   213 <- 214 (and 214 <- 213)
#   222 <- 223 (and 223 <- 222)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowData($row) {
    return $ws.Range("B$row`:AD$row").Value2
}

function Set-RowData($row, $data) {
    $ws.Range("B$row`:AD$row").Value2 = $data
}

# --- Capture all source data first, before any writes happen ---
$data112 = Get-RowData 112
$data113 = Get-RowData 113

$data124 = Get-RowData 124
$data125 = Get-RowData 125

$data143 = Get-RowData 143
$data144 = Get-RowData 144

$data175 = Get-RowData 175
$data179 = Get-RowData 179
$data180 = Get-RowData 180

$data213 = Get-RowData 213
$data214 = Get-RowData 214

$data222 = Get-RowData 222
$data223 = Get-RowData 223

# --- Apply swaps / rotation using the captured snapshots ---

# Simple pairwise swaps
Set-RowData 112 $data113
Set-RowData 113 $data112

Set-RowData 124 $data125
Set-RowData 125 $data124

Set-RowData 143 $data144
Set-RowData 144 $data143

Set-RowData 213 $data214
Set-RowData 214 $data213

Set-RowData 222 $data223
Set-RowData 223 $data222

# 3-way rotation: 175 <- 180, 179 <- 175, 180 <- 179
Set-RowData 175 $data180
Set-RowData 179 $data175
Set-RowData 180 $data179
